# Grading spreadsheet update: fill in missing "Лаба №2" (column D) scores (and a
# few corrections elsewhere) after students re-submitted / copied over their
# answers, and update the free-text notes in column M accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- Row 2: Алсуфьева Мария Николаевна ---
$ws.Range("B2").Value = 4
$ws.Range("D2").Value = 5

# --- Row 3: Ардаширова Амина Рифовна ---
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 0
$ws.Range("M3").Value = "переписаны верно все номера"

# --- Row 6: Быков Вадим Дмитриевич ---
$ws.Range("B6").Value = -1
$ws.Range("D6").Value = 5

# --- Row 7: Галямова Яна Дмитриевна ---
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 5
$ws.Range("M7").Value = "переписаны верно все номера"

# --- Row 9: Жамсаранова Аяна Жаргаловна ---
$ws.Range("B9").Value = 4
$ws.Range("D9").Value = 5

# --- Row 10: Захаренкова Екатерина Денисовна ---
$ws.Range("D10").Value = 5

# --- Row 11: Иванов Дмитрий Сергеевич ---
$ws.Range("B11").Value = 4
$ws.Range("D11").Value = 5

# --- Row 12: Косарынская Анна Александровна ---
$ws.Range("B12").Value = 4

# --- Row 13: Костылев Владимир Алексеевич ---
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 0
$ws.Range("M13").Value = "переписаны верно все номера"

# --- Row 14: Круглов Кирилл Максимович ---
$ws.Range("D14").Value = 0

# --- Row 15: Крутов Никита Сергеевич ---
$ws.Range("D15").Value = 5

# --- Row 16: Ларюшин Виктор Романович ---
$ws.Range("E16").Value = 5

# --- Row 17: Махаури Амина Эдуардовна ---
$ws.Range("D17").Value = 5

# --- Row 18: Мачкалян Тигран Норайрович ---
$ws.Range("D18").Value = 5

# --- Row 19: Молокова Татьяна Михайловна ---
$ws.Range("B19").Value = 5
$ws.Range("D19").Value = 5

# --- Row 20: Нефодина Анна Александровна ---
$ws.Range("D20").Value = 5

# --- Row 21: Ротанкова Вера Владимировна ---
$ws.Range("B21").Value = 5
$ws.Range("D21").Value = 5

# --- Row 23: Смирнова Юлия Максимовна ---
$ws.Range("D23").Value = 5

# --- Row 24: Сычиков Владимир Андреевич ---
$ws.Range("D24").Value = 5

# --- Row 25: Ушакова Александра Юрьевна ---
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 5
$ws.Range("M25").Value = "переписаны верно все номера"

# Leave the final selection on the last-edited note cell, matching the source
# workbook's saved cursor position.
$ws.Range("M13").Select() | Out-Null
